$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44189
$ws.Range("M2").Value = 40

# Row 4
$ws.Range("D4").Value = 44186
$ws.Range("M4").Value = 40

# Row 5
$ws.Range("D5").Value = 44188
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 3000

# Row 6
$ws.Range("D6").Value = 44175
$ws.Range("M6").Value = 25
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("S6").Value = 4000

# Row 7
$ws.Range("D7").Value = 44196
$ws.Range("M7").Value = 56
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 3000

# Row 8
$ws.Range("D8").Value = 44193

# Row 9
$ws.Range("D9").Value = 44179
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 4000
